# Applies the scheduled-runner data refresh to the Cerberus_Profits market
# sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Each sheet tracks a leve's
# current market price columns (H:N); this run's pull refreshed the
# average-price inputs and their downstream profit computations for a
# handful of rows per sheet. Values are written directly (source data has
# no formulas). One row (WVR!126) lost its HQ-profit column (N) in this
# pull - NQ profit (M) absorbs the updated figure instead.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 7286.143
$ws.Range("I9").Value = 10138.4
$ws.Range("J9").Value = 155.5
$ws.Range("K9").Value = 10138.4
$ws.Range("L9").Value = 155.5
$ws.Range("M9").Value = -9969.4
$ws.Range("N9").Value = -493.5
$ws.Range("I19").Value = 1098.75
$ws.Range("J19").Value = 1269.3334
$ws.Range("K19").Value = 1098.75
$ws.Range("L19").Value = 1269.3334
$ws.Range("M19").Value = -923.75
$ws.Range("N19").Value = -1619.3334
$ws.Range("H62").Value = 125001500
$ws.Range("I62").Value = 125001500
$ws.Range("K62").Value = 125001500
$ws.Range("M62").Value = -125000876
$ws.Range("H65").Value = 125001500
$ws.Range("I65").Value = 125001500
$ws.Range("K65").Value = 625007500
$ws.Range("M65").Value = -625004380
$ws.Range("H70").Value = 3813.2046
$ws.Range("I70").Value = 1387
$ws.Range("K70").Value = 4161
$ws.Range("M70").Value = -3891
$ws.Range("H73").Value = 3813.2046
$ws.Range("I73").Value = 1387
$ws.Range("K73").Value = 4161
$ws.Range("M73").Value = -3225
$ws.Range("H116").Value = 8687.75
$ws.Range("I116").Value = 24500
$ws.Range("J116").Value = 3417
$ws.Range("K116").Value = 24500
$ws.Range("L116").Value = 3417
$ws.Range("M116").Value = -21058
$ws.Range("N116").Value = -10301
$ws.Range("H134").Value = 34998.688
$ws.Range("J134").Value = 34998.688
$ws.Range("L134").Value = 34998.688
$ws.Range("N134").Value = -45138.688
$ws.Range("H138").Value = 4497.606
$ws.Range("I138").Value = 4031.4443
$ws.Range("J138").Value = 5057
$ws.Range("K138").Value = 12094.3329
$ws.Range("L138").Value = 15171
$ws.Range("M138").Value = -6954.332900000001
$ws.Range("N138").Value = -25451

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1193.2174
$ws.Range("I97").Value = 640.5294
$ws.Range("J97").Value = 2759.1667
$ws.Range("K97").Value = 640.5294
$ws.Range("L97").Value = 2759.1667
$ws.Range("M97").Value = -144.5294
$ws.Range("N97").Value = -3751.1667

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3060.8125
$ws.Range("I105").Value = 3212.5
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 3212.5
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = -1465.5
$ws.Range("N105").Value = -5493
$ws.Range("H133").Value = 87999
$ws.Range("J133").Value = 87999
$ws.Range("L133").Value = 87999
$ws.Range("N133").Value = -98119

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2129.1538
$ws.Range("I31").Value = 1615.3334
$ws.Range("J31").Value = 2829.818
$ws.Range("K31").Value = 1615.3334
$ws.Range("L31").Value = 2829.818
$ws.Range("M31").Value = -1320.3334
$ws.Range("N31").Value = -3419.818
$ws.Range("H34").Value = 2129.1538
$ws.Range("I34").Value = 1615.3334
$ws.Range("J34").Value = 2829.818
$ws.Range("K34").Value = 1615.3334
$ws.Range("L34").Value = 2829.818
$ws.Range("M34").Value = -1413.3334
$ws.Range("N34").Value = -3233.818
$ws.Range("H86").Value = 6662.857
$ws.Range("I86").Value = 6762.5
$ws.Range("J86").Value = 6530
$ws.Range("K86").Value = 6762.5
$ws.Range("L86").Value = 6530
$ws.Range("M86").Value = -5639.5
$ws.Range("N86").Value = -8776
$ws.Range("H89").Value = 6662.857
$ws.Range("I89").Value = 6762.5
$ws.Range("J89").Value = 6530
$ws.Range("K89").Value = 33812.5
$ws.Range("L89").Value = 32650
$ws.Range("M89").Value = -28196.5
$ws.Range("N89").Value = -43882
$ws.Range("H134").Value = 1915.2727
$ws.Range("I134").Value = 1856.8
$ws.Range("K134").Value = 5570.4
$ws.Range("M134").Value = -3035.4

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 391.26666
$ws.Range("J5").Value = 348.33334
$ws.Range("L5").Value = 1045.00002
$ws.Range("N5").Value = -1269.00002
$ws.Range("H87").Value = 9308.691999999999
$ws.Range("I87").Value = 6779.222
$ws.Range("K87").Value = 20337.666
$ws.Range("M87").Value = -19089.666
$ws.Range("H90").Value = 9308.691999999999
$ws.Range("I90").Value = 6779.222
$ws.Range("K90").Value = 61012.998
$ws.Range("M90").Value = -54772.998
$ws.Range("H122").Value = 246.75
$ws.Range("J122").Value = 266.33334
$ws.Range("L122").Value = 2397.00006
$ws.Range("N122").Value = -7297.00006
$ws.Range("H135").Value = 391.26666
$ws.Range("J135").Value = 348.33334
$ws.Range("L135").Value = 3135.00006
$ws.Range("N135").Value = -8205.00006

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3180
$ws.Range("J80").Value = 3999.5
$ws.Range("L80").Value = 3999.5
$ws.Range("N80").Value = -5995.5
$ws.Range("H83").Value = 3180
$ws.Range("J83").Value = 3999.5
$ws.Range("L83").Value = 19997.5
$ws.Range("N83").Value = -29981.5
$ws.Range("H102").Value = 14035.75
$ws.Range("I102").Value = 25587.637
$ws.Range("K102").Value = 25587.637
$ws.Range("M102").Value = -23965.637
$ws.Range("H132").Value = 2665.182
$ws.Range("I132").Value = 2418.8235
$ws.Range("J132").Value = 3502.8
$ws.Range("K132").Value = 7256.470499999999
$ws.Range("L132").Value = 10508.4
$ws.Range("M132").Value = -4726.470499999999
$ws.Range("N132").Value = -15568.4
$ws.Range("H135").Value = 91511.8
$ws.Range("J135").Value = 91511.8
$ws.Range("L135").Value = 91511.8
$ws.Range("N135").Value = -101651.8

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6797.5
$ws.Range("I7").Value = 6694.1665
$ws.Range("J7").Value = 6875
$ws.Range("K7").Value = 6694.1665
$ws.Range("L7").Value = 6875
$ws.Range("M7").Value = -6582.1665
$ws.Range("N7").Value = -7099
$ws.Range("H68").Value = 2356.158
$ws.Range("I68").Value = 2327.6428
$ws.Range("K68").Value = 2327.6428
$ws.Range("M68").Value = -1578.6428
$ws.Range("H71").Value = 2356.158
$ws.Range("I71").Value = 2327.6428
$ws.Range("K71").Value = 11638.214
$ws.Range("M71").Value = -7894.214
$ws.Range("H109").Value = 63420
$ws.Range("J109").Value = 63420
$ws.Range("L109").Value = 63420
$ws.Range("N109").Value = -66194
$ws.Range("H126").Value = 6797.5
$ws.Range("I126").Value = 6694.1665
$ws.Range("J126").Value = 6875
$ws.Range("K126").Value = 20082.4995
$ws.Range("L126").Value = 20625
$ws.Range("M126").Value = -17612.4995
$ws.Range("N126").Value = -25565

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6473.75
$ws.Range("I62").Value = 5758.4
$ws.Range("K62").Value = 5758.4
$ws.Range("M62").Value = -5134.4
$ws.Range("H65").Value = 6473.75
$ws.Range("I65").Value = 5758.4
$ws.Range("K65").Value = 28792
$ws.Range("M65").Value = -25672
$ws.Range("H81").Value = 3675.8823
$ws.Range("I81").Value = 4950.75
$ws.Range("K81").Value = 9901.5
$ws.Range("M81").Value = -8840.5
$ws.Range("H84").Value = 3675.8823
$ws.Range("I84").Value = 4950.75
$ws.Range("K84").Value = 49507.5
$ws.Range("M84").Value = -44203.5
$ws.Range("H109").Value = 69376.5
$ws.Range("J109").Value = 69376.5
$ws.Range("L109").Value = 69376.5
$ws.Range("N109").Value = -72150.5
$ws.Range("H126").Value = 2769.8
$ws.Range("I126").Value = 2769.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8309.400000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5839.400000000001
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 3655.625
$ws.Range("I136").Value = 3829
$ws.Range("K136").Value = 11487
$ws.Range("M136").Value = -8937
$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280
